# Import Viral Results verifies Well Plate Barcode and Position given when
# providing the other (CVDLS-201).
#
# The workbook grows from 5 data rows to 10: row 3 (old "TubeQPCRResults0002")
# splits into a blank spacer row (new row 3), a row missing its Specimen ID
# (new row 4, highlighted red) and the original row shifted down (new row 5).
# Old rows 4 and 5 shift down to rows 6 and 7. Three brand new rows (8, 9, 10)
# are appended, each missing exactly one of {Result, Position, Well Plate
# Barcode} and highlighted red in that cell.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$RED = 255

# ---------------------------------------------------------------------
# 1) Header row (row 1): bold + underline every header cell A1:M1
# ---------------------------------------------------------------------
$ws.Range("A1:M1").Font.Bold = $true
$ws.Range("A1:M1").Font.Underline = $true

# ---------------------------------------------------------------------
# 2) Shift the old rows 3,4,5 (TubeQPCRResults0002/3/4) down to 5,6,7.
#    Old row 2 (TubeQPCRResults0001) stays put. Copy from the bottom up
#    so a source row is never clobbered before it has been read.
# ---------------------------------------------------------------------
$ws.Range("A5:M5").Copy()
$ws.Range("A7:M7").PasteSpecial(-4104)

$ws.Range("A4:M4").Copy()
$ws.Range("A6:M6").PasteSpecial(-4104)

$ws.Range("A3:M3").Copy()
$ws.Range("A5:M5").PasteSpecial(-4104)
$excel.CutCopyMode = 0

# ---------------------------------------------------------------------
# 3) New row 4: same data as row 5 (TubeQPCRResults0002, still intact at
#    this point) but Specimen ID (A4) is blank/red, Source Well becomes
#    "B1" (a new value), and the yellow highlight from the original row
#    is removed (plain white).
# ---------------------------------------------------------------------
$ws.Range("A3:M3").Copy()
$ws.Range("A4:M4").PasteSpecial(-4104)
$excel.CutCopyMode = 0
$ws.Range("A4").ClearContents()
$ws.Range("A4:M4").Interior.Pattern = -4142
$ws.Range("A4").NumberFormat = "0"
$ws.Range("A4").Interior.Color = $RED
$ws.Range("C4").Value = "B1"

# ---------------------------------------------------------------------
# 4) New row 3: blank spacer row - clear content/format so only A3 keeps
#    an (invisible) numeric-format style, matching the original row 3's
#    number format.
# ---------------------------------------------------------------------
$ws.Range("A3:M3").ClearContents()
$ws.Range("B3:M3").ClearFormats()
$ws.Range("A3").NumberFormat = "0"

# ---------------------------------------------------------------------
# 5) New rows 8, 9, 10 - append new specimens.
# ---------------------------------------------------------------------

# Row 8: TubeQPCRResults0005 - missing Result (B8, red); Position=F8
$ws.Range("A8").Value = "TubeQPCRResults0005"
$ws.Range("B8").Value = ""
$ws.Range("C8").Value = "F8"
$ws.Range("D8").Value = "20200715_PVE690RLR_1"
$ws.Range("E8").Value = "QPCRResults"
$ws.Range("F8").Value = "Rack001"
$ws.Range("G8").Value = "B02"
$ws.Range("H8").Value = "Undetermined"
$ws.Range("I8").Value = 500
$ws.Range("J8").Value = 21.987654320000001
$ws.Range("K8").Value = 0
$ws.Range("L8").Value = 600
$ws.Range("M8").Value = 0
$ws.Range("A8").NumberFormat = "0"
$ws.Range("B8").Interior.Color = $RED

# Row 9: TubeQPCRResults0006 - Barcode given (E9), Position missing (C9, red)
$ws.Range("A9").Value = "TubeQPCRResults0006"
$ws.Range("B9").Value = "Positive"
$ws.Range("C9").Value = ""
$ws.Range("D9").Value = "20200715_PVE690RLR_1"
$ws.Range("E9").Value = "QPCRResults"
$ws.Range("F9").Value = "Rack001"
$ws.Range("G9").Value = "B02"
$ws.Range("H9").Value = "Undetermined"
$ws.Range("I9").Value = 500
$ws.Range("J9").Value = 21.987654320000001
$ws.Range("K9").Value = 0
$ws.Range("L9").Value = 600
$ws.Range("M9").Value = 0
$ws.Range("A9").NumberFormat = "0"
$ws.Range("C9").Interior.Color = $RED

# Row 10: TubeQPCRResults0007 - Position given (C10), Barcode missing (E10, red)
$ws.Range("A10").Value = "TubeQPCRResults0007"
$ws.Range("B10").Value = "Positive"
$ws.Range("C10").Value = "H10"
$ws.Range("D10").Value = "20200715_PVE690RLR_1"
$ws.Range("E10").Value = ""
$ws.Range("F10").Value = "Rack001"
$ws.Range("G10").Value = "B02"
$ws.Range("H10").Value = "Undetermined"
$ws.Range("I10").Value = 500
$ws.Range("J10").Value = 21.987654320000001
$ws.Range("K10").Value = 0
$ws.Range("L10").Value = 600
$ws.Range("M10").Value = 0
$ws.Range("A10").NumberFormat = "0"
$ws.Range("E10").Interior.Color = $RED

# ---------------------------------------------------------------------
# 6) Selection lands on C9 (the newly-flagged missing Position cell).
# ---------------------------------------------------------------------
$ws.Range("C9").Select()
